# The workbook's single worksheet is being renamed from "Sheet1" to
# "testunsafe" to reflect the contents of the data file
# (data/erp/testunsafedata.xlsx).
$wb = $excel.ActiveWorkbook

$ws = $wb.ActiveSheet
if ($wb.Worksheets | Where-Object { $_.Name -eq "Sheet1" }) {
    $ws = $wb.Worksheets.Item("Sheet1")
}

$ws.Name = "testunsafe"
